$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Orders")
$ws.Activate()

# "Gula Malacca" (row 10) was mis-categorized as a Seasoning; re-classify it
# as a Beverage before the bulk rename below.
$ws.Range("C10").Value = "Beverage"

# Rename the "Beverage" product-type label to "Beverages" everywhere it
# appears in the Orders sheet (the Type column, C2:C25).
$ws.Cells.Replace("Beverage", "Beverages")

# Leave the selection on the last cell touched by the replace (C20 is the
# final row whose Type was "Beverage").
$ws.Range("C20").Select()
